$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 8 (years 2003-2009), shifting remaining rows (2010-2012) up
$ws.Range("A2:A8").EntireRow.Delete()

# Update the dimension will be handled automatically by Excel on save
